# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp banner (A1) ---
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 19 de Mayo de 2020 a las 19:35"

# --- Reorder countries: Birmania now listed before Martinica ---
# (row 153 was Martinica / row 154 was Birmania -> swap display order)
$ws.Cells.Item(153, 1).Value2 = "Birmania"
$ws.Cells.Item(154, 1).Value2 = "Martinica"

# --- Update covid numbers (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) per row ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value2 = 1558343
$ws.Cells.Item(4, 3).Value2 = 8049
$ws.Cells.Item(4, 5).Value2 = 1105800
$ws.Cells.Item(4, 7).Value2 = 504
$ws.Cells.Item(4, 8).Value2 = 92485

# Row 7: Brasil
$ws.Cells.Item(7, 2).Value2 = 262545
$ws.Cells.Item(7, 3).Value2 = 7177
$ws.Cells.Item(7, 5).Value2 = 144577
$ws.Cells.Item(7, 7).Value2 = 656
$ws.Cells.Item(7, 8).Value2 = 17509

# Row 8: Reino Unido
$ws.Cells.Item(8, 2).Value2 = 248818
$ws.Cells.Item(8, 3).Value2 = 2412
$ws.Cells.Item(8, 7).Value2 = 545
$ws.Cells.Item(8, 8).Value2 = 35341

# Row 11: Alemania
$ws.Cells.Item(11, 2).Value2 = 177696
$ws.Cells.Item(11, 3).Value2 = 407
$ws.Cells.Item(11, 5).Value2 = 13835
$ws.Cells.Item(11, 7).Value2 = 38
$ws.Cells.Item(11, 8).Value2 = 8161

# Row 14: India
$ws.Cells.Item(14, 2).Value2 = 103935
$ws.Cells.Item(14, 3).Value2 = 3607
$ws.Cells.Item(14, 4).Value2 = 40865
$ws.Cells.Item(14, 5).Value2 = 59854
$ws.Cells.Item(14, 7).Value2 = 60
$ws.Cells.Item(14, 8).Value2 = 3216

# Row 25: Ecuador
$ws.Cells.Item(25, 2).Value2 = 34151
$ws.Cells.Item(25, 3).Value2 = 569
$ws.Cells.Item(25, 4).Value2 = 3457
$ws.Cells.Item(25, 5).Value2 = 27855
$ws.Cells.Item(25, 7).Value2 = 40
$ws.Cells.Item(25, 8).Value2 = 2839

# Row 33: Irlanda
$ws.Cells.Item(33, 2).Value2 = 24251
$ws.Cells.Item(33, 3).Value2 = 51
$ws.Cells.Item(33, 5).Value2 = 3220
$ws.Cells.Item(33, 7).Value2 = 14
$ws.Cells.Item(33, 8).Value2 = 1561

# Row 51: Chequia
$ws.Cells.Item(51, 2).Value2 = 8630
$ws.Cells.Item(51, 3).Value2 = 44
$ws.Cells.Item(51, 4).Value2 = 5725
$ws.Cells.Item(51, 5).Value2 = 2603
$ws.Cells.Item(51, 7).Value2 = 5
$ws.Cells.Item(51, 8).Value2 = 302

# Row 86: Cuba
$ws.Cells.Item(86, 2).Value2 = 1887
$ws.Cells.Item(86, 3).Value2 = 6
$ws.Cells.Item(86, 4).Value2 = 1538
$ws.Cells.Item(86, 5).Value2 = 270

# Row 105: Sri Lanka
$ws.Cells.Item(105, 2).Value2 = 1021
$ws.Cells.Item(105, 3).Value2 = 29
$ws.Cells.Item(105, 5).Value2 = 443

# Row 153: now Birmania (updated figures)
$ws.Cells.Item(153, 2).Value2 = 193
$ws.Cells.Item(153, 3).Value2 = 5
$ws.Cells.Item(153, 4).Value2 = 104
$ws.Cells.Item(153, 5).Value2 = 83
$ws.Cells.Item(153, 8).Value2 = 6

# Row 154: now Martinica (figures carried over from prior Martinica row)
$ws.Cells.Item(154, 2).Value2 = 192
$ws.Cells.Item(154, 3).Value2 = 0
$ws.Cells.Item(154, 4).Value2 = 91
$ws.Cells.Item(154, 5).Value2 = 87
$ws.Cells.Item(154, 8).Value2 = 14
